# Gatchina city 17-19 (added)
# - highlight G56 (the "livarea" value for 2014 on the Гатчинский МР row) in green
# - add three new rows (57, 58, 59) of data for "Гатчина" (2019, 2018, 2017)
# - move the active selection/view near the new rows

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Highlight G56 with a green fill (adds a new fill + cellXf) ---
# Alignment is (re-)applied first so the engine reuses the already-centered
# style, then Interior.Color mints exactly one new fill + cellXf for it.
$ws.Range("G56").HorizontalAlignment = -4108  # xlCenter
$ws.Range("G56").Interior.Color = 5296274     # RGB(146, 208, 80) == 0xFF92D050

# --- Fill in the "???" placeholder cells that line up with the rest of the table for row 56 ---
$ws.Range("H56:K56").HorizontalAlignment = -4108
$ws.Range("H56").Value = "???"
$ws.Range("I56").Value = "???"
$ws.Range("J56").Value = "???"
$ws.Range("K56").Value = "???"

# L56 picks up the shared "centered" style used by the rest of the numeric columns
$ws.Range("L56").HorizontalAlignment = -4108  # xlCenter

# --- New rows of data for "Гатчина" ---
$gatchinaRows = @(
    @{ Row=57; Year=2019; PopSize=91.677999999999997; AvgEmployers=27.295000000000002; Unemployed=189; AvgSalary=45526; Invests=3022.6; Companies=5566; FactoriesCap=22280.27; RetailFormula="=21358889.6/1000"; FoodservFormula="=333725.2/1000"; Saldo=-1280 },
    @{ Row=58; Year=2018; PopSize=93.721999999999994; AvgEmployers=27.602; Unemployed=144; AvgSalary=41230.699999999997; Invests=12843.2; Companies=5245; FactoriesCap=33385.4; RetailFormula="=15613892/1000"; FoodservFormula="=170054/1000"; Saldo=-28 },
    @{ Row=59; Year=2017; PopSize=94.45; AvgEmployers=22.256; Unemployed=1006; AvgSalary=41319.300000000003; Invests=2983.3; Companies=5071; FactoriesCap=26147.7; RetailFormula="= 9892900 / 1000"; FoodservFormula="= 130235 / 1000"; Saldo=-114 }
)

foreach ($row in $gatchinaRows) {
    $r = $row.Row.ToString()

    # Centered alignment for the whole row first, so plain numeric cells reuse
    # the existing "style 2" xf instead of minting new ones.
    $ws.Range("A" + $r + ":U" + $r).HorizontalAlignment = -4108
    # A/B (name + whole year) reuse the existing bold "style 1" xf.
    $ws.Range("A" + $r + ":B" + $r).Font.Bold = $true

    $ws.Range("A" + $r).Value = "Гатчина"
    $ws.Range("B" + $r).Value = $row.Year
    $ws.Range("C" + $r).Value = $row.PopSize
    $ws.Range("D" + $r).Value = $row.AvgEmployers
    $ws.Range("E" + $r).Value = $row.Unemployed
    $ws.Range("F" + $r).Value = $row.AvgSalary
    $ws.Range("G" + $r).Value = "???"
    $ws.Range("H" + $r).Value = "???"
    $ws.Range("I" + $r).Value = "???"
    $ws.Range("J" + $r).Value = "???"
    $ws.Range("K" + $r).Value = "???"
    $ws.Range("L" + $r).Value = $row.Invests
    $ws.Range("M" + $r).Value = "???"
    $ws.Range("N" + $r).Value = $row.Companies
    $ws.Range("O" + $r).Value = $row.FactoriesCap
    $ws.Range("P" + $r).Value = "???"
    $ws.Range("Q" + $r).Value = "???"
    $ws.Range("R" + $r).Value = "???"
    $ws.Range("S" + $r).Formula = $row.RetailFormula
    $ws.Range("T" + $r).Formula = $row.FoodservFormula
    $ws.Range("U" + $r).Value = $row.Saldo

    # N, P, Q, R reuse the existing yellow-highlighted centered style (index 3).
    $ws.Range("N" + $r).Interior.ColorIndex = 6
    $ws.Range("P" + $r + ":R" + $r).Interior.ColorIndex = 6
}

# --- Selection / view ---
$ws.Range("K65").Select()

"Gatchina rows 57-59 added"
